$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("리퍼리움", "엠블", "rare"),
    @("스톰엑스", "아이오에스티", "rare"),
    @("엘프", "골렘", "rare"),
    @("아더", "아르고", "rare"),
    @("스팀", "하이브", "rare"),
    @("메탈", "아크", "rare"),
    @("어거", "웨이브", "rare"),
    @("코박토큰", "알파쿼크", "rare"),
    @("플로우", "던프로토콜", "rare"),
    @("샌드박스", "디센트럴랜드", "rare"),
    @("세럼", "카바", "rare"),
    @("퀀텀", "네오", "rare"),
    @("트론", "비트토렌트", "rare"),
    @("메디블록", "메타디움", "rare")
)

$row = 5
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$ws.Range("C8").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

